$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '62.896.56'
Set-TextValue 'E2' '  +5.08%  '
Set-TextValue 'D3' '3.351.93'
Set-TextValue 'E3' '  +5.07%  '
Set-TextValue 'E4' '  +0.03%  '
Set-TextValue 'D5' '556.97'
Set-TextValue 'E5' '  +3.70%  '
Set-TextValue 'D6' '153.06'
Set-TextValue 'E6' '  +5.91%  '
Set-TextValue 'D7' '0.999'
Set-TextValue 'E7' '  -0.16%  '
Set-TextValue 'D8' '0.529'
Set-TextValue 'E8' '  +0.59%  '
Set-TextValue 'E9' '  +2.34%  '
Set-TextValue 'E10' '  +4.77%  '
Set-TextValue 'D11' '0.441'
Set-TextValue 'E11' '  +2.19%  '
Set-TextValue 'D12' '3.933.48'
Set-TextValue 'E12' '  +5.05%  '
Set-TextValue 'E13' '  +0.32%  '
Set-TextValue 'D14' '0.0000181'
Set-TextValue 'E14' '  +4.12%  '
Set-TextValue 'D15' '27.01'
Set-TextValue 'E15' '  +3.72%  '
Set-TextValue 'D16' '63.002.94'
Set-TextValue 'E16' '  +5.15%  '
Set-TextValue 'D17' '3.356.37'
Set-TextValue 'E17' '  +4.69%  '
Set-TextValue 'D18' '6.50'
Set-TextValue 'E18' '  +4.26%  '
Set-TextValue 'E19' '  +5.25%  '
Set-TextValue 'D20' '8.45'
Set-TextValue 'E20' '  +0.94%  '
Set-TextValue 'D21' '389.15'
Set-TextValue 'E21' '  +1.69%  '
Set-TextValue 'B22' 'Polygon'
Set-TextValue 'C22' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D22' '0.542'
Set-TextValue 'E22' '  +2.52%  '
Set-TextValue 'B23' 'Dai'
Set-TextValue 'C23' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D23' '0.999'
Set-TextValue 'E23' '  -0.13%  '
Set-TextValue 'D24' '70.69'
Set-TextValue 'E24' '  +0.59%  '
Set-TextValue 'D25' '0.180'
Set-TextValue 'E25' '  +5.35%  '
Set-TextValue 'D26' '8.83'
Set-TextValue 'E26' '  -0.27%  '
Set-TextValue 'D27' '0.0₃0970'
Set-TextValue 'E27' '  +8.15%  '
Set-TextValue 'E28' '  +0.18%  '
Set-TextValue 'E29' '  +4.60%  '
Set-TextValue 'D30' '6.48'
Set-TextValue 'E30' '  +5.48%  '
Set-TextValue 'D31' '5.64'
Set-TextValue 'E31' '  +4.92%  '
Set-TextValue 'D32' '23.06'
Set-TextValue 'E32' '  +3.01%  '
Set-TextValue 'D33' '1.31'
Set-TextValue 'E33' '  +8.18%  '
Set-TextValue 'E34' '  +2.84%  '
Set-TextValue 'E35' '  +10.29%  '
Set-TextValue 'D36' '159.63'
Set-TextValue 'E36' '  +2.23%  '
Set-TextValue 'E37' '  +12.23%  '
Set-TextValue 'D38' '27.32'
Set-TextValue 'E38' '  +6.81%  '
Set-TextValue 'E39' '  +4.58%  '
Set-TextValue 'D40' '2.846.14'
Set-TextValue 'E40' '  +2.77%  '
Set-TextValue 'E41' '  +8.55%  '
Set-TextValue 'D42' '4.33'
Set-TextValue 'E42' '  +1.47%  '
Set-TextValue 'D43' '40.73'
Set-TextValue 'E43' '  +2.44%  '
Set-TextValue 'D44' '0.748'
Set-TextValue 'E44' '  +2.68%  '
Set-TextValue 'E45' '  +3.98%  '
Set-TextValue 'B46' 'InjectiveProtocol'
Set-TextValue 'C46' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D46' '22.12'
Set-TextValue 'E46' '  +7.98%  '
Set-TextValue 'B47' 'RenzoRestakedETH'
Set-TextValue 'C47' 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue 'D47' '3.399.09'
Set-TextValue 'E47' '  +5.15%  '
Set-TextValue 'E48' '  +2.40%  '
Set-TextValue 'D49' '6.30'
Set-TextValue 'E49' '  +1.98%  '
Set-TextValue 'D50' '0.811'
Set-TextValue 'E50' '  +1.60%  '
Set-TextValue 'D51' '281.09'
Set-TextValue 'E51' '  +7.28%  '
